# GDE-6731 Update to takeout screenshot keywords in 04DealNotebook.robot
#
# The deal / facility identifiers that are re-generated each UAT run are
# bumped to the latest run's values, and the customer "ARIHANT TRADING
# COMPANY 000008" is replaced by the newer "ARIHANT TRADING COMPANY 000010"
# record across every sheet that references it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ORIG03_Customer
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ORIG03_Customer")
$ws.Range("C2").Value = "ARIHANT TRADING COMPANY 000010"
$ws.Range("E2").Value = 70
$ws.Range("O2").Value = "ARIHANT TRADING COMPANY 000010"
$ws.Range("AJ2").Value = "Foreign Government, Foreign Government"

# ---------------------------------------------------------------------
# CRED01_DealSetup
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRED01_DealSetup")
$ws.Range("E2").Value = "UAT4_25082020161406"
$ws.Range("F2").Value = "UAT425082020161406"
$ws.Range("G2").Value = "FAC-A_25082020162140"
$ws.Range("J2").Value = "ARIHANT TRADING COMPANY 000010"
$ws.Range("G3").Value = "FAC-B_25082020162950"

# ---------------------------------------------------------------------
# CRED02_FacilitySetup
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRED02_FacilitySetup")
$ws.Range("D2").Value = "UAT4_25082020161406"
$ws.Range("E2").Value = "FAC-A_25082020162140"
$ws.Range("Q2").Value = "ARIHANT TRADING COMPANY 000010"
$ws.Range("D3").Value = "UAT4_25082020161406"
$ws.Range("E3").Value = "FAC-B_25082020162950"
$ws.Range("Q3").Value = "ARIHANT TRADING COMPANY 000010"

# ---------------------------------------------------------------------
# SERV01_LoanDrawdown
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SERV01_LoanDrawdown")
$ws.Range("C2").Value = "UAT4_25082020161406"
$ws.Range("D2").Value = "FAC-A_25082020162140"
$ws.Range("F2").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C3").Value = "UAT4_25082020161406"
$ws.Range("D3").Value = "FAC-B_25082020162950"
$ws.Range("F3").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C4").Value = "UAT4_25082020161406"
$ws.Range("D4").Value = "FAC-B_25082020162950"
$ws.Range("F4").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C5").Value = "UAT4_25082020161406"
$ws.Range("D5").Value = "FAC-B_25082020162950"
$ws.Range("F5").Value = "ARIHANT TRADING COMPANY 000010"

# ---------------------------------------------------------------------
# UAT04_Fees
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("UAT04_Fees")
$ws.Range("C2").Value = "UAT4_25082020161406"
$ws.Range("D2").Value = "FAC-A_25082020162140"
$ws.Range("E2").Value = "ARIHANT TRADING COMPANY 000010"
$ws.Range("F2").Value = "ONG000000000324"

$ws.Range("C3").Value = "UAT4_25082020161406"
$ws.Range("D3").Value = "FAC-B_25082020162950"
$ws.Range("E3").Value = "ARIHANT TRADING COMPANY 000010"
$ws.Range("F3").Value = "ONG000000000325"

# ---------------------------------------------------------------------
# UAT04_Runbook
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("UAT04_Runbook")

$ws.Range("C2").Value = "UAT4_25082020161406"
$ws.Range("D2").Value = "FAC-A_25082020162140"
$ws.Range("G2").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C3").Value = "UAT4_25082020161406"
$ws.Range("D3").Value = "FAC-A_25082020162140"
$ws.Range("G3").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C4").Value = "UAT4_25082020161406"
$ws.Range("D4").Value = "FAC-A_25082020162140"
$ws.Range("G4").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C5").Value = "UAT4_25082020161406"
$ws.Range("D5").Value = "FAC-A_25082020162140"
$ws.Range("G5").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C6").Value = "UAT4_25082020161406"
$ws.Range("D6").Value = "FAC-B_25082020162950"
$ws.Range("G6").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C7").Value = "UAT4_25082020161406"
$ws.Range("D7").Value = "FAC-B_25082020162950"
$ws.Range("G7").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C8").Value = "UAT4_25082020161406"
$ws.Range("D8").Value = "FAC-B_25082020162950"
$ws.Range("G8").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C9").Value = "UAT4_25082020161406"
$ws.Range("D9").Value = "FAC-B_25082020162950"
$ws.Range("G9").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C10").Value = "UAT4_25082020161406"
$ws.Range("D10").Value = "FAC-B_25082020162950"
$ws.Range("G10").Value = "ARIHANT TRADING COMPANY 000010"

$ws.Range("C11").Value = "UAT4_25082020161406"
$ws.Range("D11").Value = "FAC-B_25082020162950"
$ws.Range("G11").Value = "ARIHANT TRADING COMPANY 000010"
